$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $origStyle = $Cell.Style
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "43.008.78"
Set-TextValue $ws.Range("E2") "  -0.49%  "

Set-TextValue $ws.Range("D3") "2.239.64"
Set-TextValue $ws.Range("E3") "  -1.74%  "

Set-TextValue $ws.Range("E4") "  -0.01%  "

Set-TextValue $ws.Range("D5") "114.61"
Set-TextValue $ws.Range("E5") "  +1.35%  "

Set-TextValue $ws.Range("D6") "269.39"
Set-TextValue $ws.Range("E6") "  +1.61%  "

Set-TextValue $ws.Range("D7") "0.629"
Set-TextValue $ws.Range("E7") "  +2.03%  "

Set-TextValue $ws.Range("E8") "  +0.35%  "

Set-TextValue $ws.Range("D9") "0.607"
Set-TextValue $ws.Range("E9") "  +0.25%  "

Set-TextValue $ws.Range("D10") "46.34"
Set-TextValue $ws.Range("E10") "  -2.44%  "

Set-TextValue $ws.Range("E11") "  -0.36%  "

Set-TextValue $ws.Range("D12") "9.17"
Set-TextValue $ws.Range("E12") "  -1.10%  "

Set-TextValue $ws.Range("E13") "  -2.54%  "

Set-TextValue $ws.Range("D14") "15.37"
Set-TextValue $ws.Range("E14") "  -0.56%  "

Set-TextValue $ws.Range("D15") "0.878"
Set-TextValue $ws.Range("E15") "  +1.74%  "

Set-TextValue $ws.Range("D16") "2.573.79"
Set-TextValue $ws.Range("E16") "  -1.31%  "

Set-TextValue $ws.Range("D17") "2.243.46"
Set-TextValue $ws.Range("E17") "  -1.32%  "

Set-TextValue $ws.Range("D18") "42.996.28"
Set-TextValue $ws.Range("E18") "  -0.46%  "

Set-TextValue $ws.Range("E19") "  -0.83%  "

Set-TextValue $ws.Range("D20") "6.76"
Set-TextValue $ws.Range("E20") "  -0.57%  "

Set-TextValue $ws.Range("D21") "72.00"
Set-TextValue $ws.Range("E21") "  +0.17%  "

Set-TextValue $ws.Range("E22") "  -4.47%  "

Set-TextValue $ws.Range("D23") "233.86"
Set-TextValue $ws.Range("E23") "  +0.07%  "

Set-TextValue $ws.Range("D24") "2.93"
Set-TextValue $ws.Range("E24") "  +1.71%  "

Set-TextValue $ws.Range("D25") "9.35"
Set-TextValue $ws.Range("E25") "  -2.54%  "

Set-TextValue $ws.Range("D26") "12.19"
Set-TextValue $ws.Range("E26") "  +6.90%  "

Set-TextValue $ws.Range("E27") "  -0.20%  "

Set-TextValue $ws.Range("D28") "40.81"
Set-TextValue $ws.Range("E28") "  -0.22%  "

Set-TextValue $ws.Range("D29") "2.22"
Set-TextValue $ws.Range("E29") "  -1.16%  "

Set-TextValue $ws.Range("E30") "  -1.87%  "

Set-TextValue $ws.Range("D31") "173.93"
Set-TextValue $ws.Range("E31") "  +0.20%  "

Set-TextValue $ws.Range("D32") "21.17"
Set-TextValue $ws.Range("E32") "  -1.49%  "

Set-TextValue $ws.Range("D33") "0.0907"
Set-TextValue $ws.Range("E33") "  +0.19%  "

Set-TextValue $ws.Range("E34") "  -2.00%  "

Set-TextValue $ws.Range("D35") "4.30"
Set-TextValue $ws.Range("E35") "  +10.56%  "

Set-TextValue $ws.Range("E36") "  +0.56%  "

Set-TextValue $ws.Range("D37") "4.69"
Set-TextValue $ws.Range("E37") "  +0.56%  "

Set-TextValue $ws.Range("D38") "0.0374"
Set-TextValue $ws.Range("E38") "  +2.30%  "

Set-TextValue $ws.Range("D39") "0.107"
Set-TextValue $ws.Range("E39") "  +2.94%  "

Set-TextValue $ws.Range("D40") "2.55"
Set-TextValue $ws.Range("E40") "  -2.78%  "

Set-TextValue $ws.Range("D41") "71.53"
Set-TextValue $ws.Range("E41") "  -6.20%  "

Set-TextValue $ws.Range("D42") "13.26"
Set-TextValue $ws.Range("E42") "  -6.21%  "

Set-TextValue $ws.Range("E43") "  -2.22%  "

Set-TextValue $ws.Range("E44") "  +0.33%  "

Set-TextValue $ws.Range("E45") "  -8.69%  "

Set-TextValue $ws.Range("D46") "1.34"

Set-TextValue $ws.Range("D47") "1.26"
Set-TextValue $ws.Range("E47") "  +0.47%  "

Set-TextValue $ws.Range("D48") "8.48"
Set-TextValue $ws.Range("E48") "  -1.77%  "

Set-TextValue $ws.Range("B49") "Cronos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0990"
Set-TextValue $ws.Range("E49") "  -0.72%  "

Set-TextValue $ws.Range("B50") "TheSandbox"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D50") "0.648"
Set-TextValue $ws.Range("E50") "  +7.66%  "

Set-TextValue $ws.Range("D51") "100.37"
Set-TextValue $ws.Range("E51") "  -3.26%  "
